$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("short term")

# Row 25: new task 93, no person
$ws.Range("A25").Value = "93. gérer le souci des décimales dans le sensitivityplot (taille de points) à partir d'ex. Ellis (pas de souci vu par Aurélie)"
$ws.Range("A25").Interior.Color = 65535
$ws.Range("A25").Borders.LineStyle = 1
$ws.Range("A25").WrapText = $true
$ws.Range("B25").Interior.Color = 65535

# Row 26: new task 94, person "ML et A"
$ws.Range("A26").Value = "94. dans les curvesplot en option ajouter un point là où la BMD est atteinte "
$ws.Range("A26").Interior.Color = 65535
$ws.Range("A26").Borders.LineStyle = 1
$ws.Range("A26").WrapText = $true
$ws.Range("B26").Value = "ML et A"
$ws.Range("B26").Interior.Color = 65535

# Row 27: new task 95, person "ML et A"
$ws.Range("A27").Value = "95. visualisation optionnelle par lignes verticales des doses testées (tous les graphes avec la dose en x)"
$ws.Range("A27").Interior.Color = 65535
$ws.Range("A27").Borders.LineStyle = 1
$ws.Range("A27").WrapText = $true
$ws.Range("B27").Value = "ML et A"
$ws.Range("B27").Interior.Color = 65535

# Row 28: new task 96, no person
$ws.Range("A28").Value = "96. ajouter dans vignette ou FAQ ex.d'utilisation de l'option ""median.and.IQR"""
$ws.Range("A28").Interior.Color = 65535
$ws.Range("A28").Borders.LineStyle = 1
$ws.Range("A28").WrapText = $true

$ws.Range("A28").Select() | Out-Null
